$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Estado" column (C) values as described by the diff:
# C2: Prestado -> Disponible
# C4: Disponible -> Prestado
# C6: Prestado -> Disponible
# C7: Disponible -> Prestado
# C9: Disponible -> Prestado
# C10: Prestado -> Disponible

$ws.Range("C2").Value = "Disponible"
$ws.Range("C4").Value = "Prestado"
$ws.Range("C6").Value = "Disponible"
$ws.Range("C7").Value = "Prestado"
$ws.Range("C9").Value = "Prestado"
$ws.Range("C10").Value = "Disponible"
